{"js": "// Update the benchmark results table: the first column of numbers (rows\n// 0-11) received refreshed measurements, and the last three rows (which\n// previously held full tab-separated result lines) are collapsed down to\n// the single summary values that used to live in rows 0-2.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// rowIndex -> new cell text\nconst changes = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"811\",\n  4: \"0.00003\",\n  5: \"0.00013\",\n  6: \"0.00006\",\n  8: \"0.00005\",\n  10: \"0.00010\",\n  11: \"0.04322\",\n  43: \"99.99\",\n  44: \"0.04\",\n  45: \"553\",\n};\n\n// Load every cell we need to touch so we can reach its body range.\nconst targetCells = [];\nfor (const idxStr of Object.keys(changes)) {\n  const idx = Number(idxStr);\n  const row = rows.items[idx];\n  row.cells.load(\"items\");\n  targetCells.push({ idx, row });\n}\nawait context.sync();\n\nfor (const { idx, row } of targetCells) {\n  const cell = row.cells.items[0];\n  const range = cell.body.getRange();\n  range.insertText(changes[idx], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the benchmark results table: the first column of numbers (rows\n# 1-12, 1-based) received refreshed measurements, and the last three rows\n# (which previously held full tab-separated result lines) are collapsed\n# down to the single summary values that used to live in rows 1-3.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# 1-based row number -> new cell text\n$changes = [ordered]@{\n  1  = \"0M\"\n  2  = \"0M\"\n  3  = \"0M\"\n  4  = \"811\"\n  5  = \"0.00003\"\n  6  = \"0.00013\"\n  7  = \"0.00006\"\n  9  = \"0.00005\"\n  11 = \"0.00010\"\n  12 = \"0.04322\"\n  44 = \"99.99\"\n  45 = \"0.04\"\n  46 = \"553\"\n}\n\nforeach ($rowNum in $changes.Keys) {\n  $cell = $tbl.Cell($rowNum, 1)\n  $cell.Range.Text = $changes[$rowNum]\n}\n"}
